$d = $word.ActiveDocument
$p1 = $d.Paragraphs.Item(1)
$r1 = $p1.Range
$r1.Collapse(0)
$r1.InsertParagraphAfter()

$p2 = $d.Paragraphs.Item(2)
$r2 = $p2.Range
$r2.Text = "No testing"

$r2b = $p2.Range
$r2b.Collapse(0)
$r2b.InsertParagraphAfter()

$p3 = $d.Paragraphs.Item(3)

# Insert the bookmark into a genuinely empty paragraph (no runs yet)
$rBookmark = $p3.Range
$rBookmark.Collapse(0)
$d.Bookmarks.Add("_GoBack", $rBookmark)

# Now set text using InsertBefore (so text gets placed ahead of bookmarkEnd, inside the run)
$rText = $d.Paragraphs.Item(3).Range
$rText.InsertBefore("Esto lo ve el profe?")
